# T1109 Contact test data - Standard User name change (6th Mar 2024)
# Replaces the Users sheet's standard-user name "Drew Koecher" with
# "Ayati Arvind", and updates the saved sheet-selection / active-tab
# state to match (Users tab active, Contact!D13 and Users!D3 selected).

$wb = $excel.ActiveWorkbook

$contactSheet = $wb.Worksheets.Item("Contact")
$usersSheet   = $wb.Worksheets.Item("Users")

# Update the Standard User name used for the external-contact test.
$usersSheet.Range("A2").Value = "Ayati Arvind"

# Restore the Contact sheet's selection (D13) while it is still active,
# then move the active tab / selection to the Users sheet (D3) - this
# mirrors the saved view state captured in the workbook after editing.
$contactSheet.Activate()
$contactSheet.Range("D13").Select()

$usersSheet.Activate()
$usersSheet.Range("D3").Select()
